# Report2.docx edit:
#   "Table with Average review r[_GoBack]atings by state"
#   "[gramStart]An[gramEnd] plot with information about detailed review ratings by state"
# becomes
#   "Table with Average review ratings by state"
#   "A plot with informati[_GoBack]on about detailed review ratings by state"
#
# i.e. the stray paragraph break inside "Average review r|atings by state" is
# healed, the grammar typo "An plot" -> "A plot" is fixed (dropping the
# now-unneeded proofErr markers around "An"), and Word's "last edit"
# _GoBack bookmark is relocated to sit at the new edit point
# ("informati|on about").

$d = $word.ActiveDocument

$full = $d.Content.Text

$oldParaA        = "Table with Average review ratings by state"
$oldParaBPrefix  = "An plot with information about "
$oldDetailed     = "detailed review ratings by state"

$idxTable = $full.IndexOf("Table with Average review r")
if ($idxTable -lt 0) {
    throw "Could not locate 'Table with Average review r' anchor text"
}

$idxDetailed = $full.IndexOf($oldDetailed, $idxTable)
if ($idxDetailed -lt 0) {
    throw "Could not locate 'detailed review ratings by state' anchor text"
}

# Confirm the span we are about to rewrite is exactly what we expect
# (paragraph A's text, the paragraph mark, then paragraph B's prefix).
$between = $d.Range($idxTable, $idxDetailed).Text
$expected = $oldParaA + "`r" + $oldParaBPrefix
if ($between -ne $expected) {
    throw "Unexpected document content before edit: [$between]"
}

# --- Step 1: rewrite paragraph A + the start of paragraph B in one go. ---
# Doing this as a single Range.Text assignment merges the runs inside each
# paragraph, drops the old (now mid-word) _GoBack bookmark and the orphaned
# proofErr gramStart/gramEnd markers around "An", while leaving the
# following "detailed review ratings by state" run (and its formatting)
# completely untouched.
$newParaA       = "Table with Average review ratings by state"
$newParaBPrefix = "A plot with information about "

$r = $d.Range($idxTable, $idxDetailed)
$r.Text = $newParaA + "`r" + $newParaBPrefix

# --- Step 2: split "A plot with informati" / "on about " into distinct
# runs by briefly dropping a bookmark right after the "A" (this forces a
# run boundary there even though both sides share identical formatting). ---
$full2 = $d.Content.Text
$idxParaB = $full2.IndexOf($newParaBPrefix, $idxTable)
if ($idxParaB -lt 0) {
    throw "Could not relocate rewritten paragraph B text"
}

$splitPoint = $idxParaB + 1   # just after the new "A"
$tmpRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("zzzTempSplit", $tmpRange) | Out-Null
$d.Bookmarks("zzzTempSplit").Delete()

# --- Step 3: drop the real _GoBack bookmark at the new edit point, between
# "informati" and "on about ". ---
$bmPoint = $idxParaB + ("A plot with informati").Length
$bmRange = $d.Range($bmPoint, $bmPoint)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Host "Updated paragraph text:" $d.Range($idxTable, $idxDetailed + $oldDetailed.Length).Text
